$d = $word.ActiveDocument
$d.Sections(1).Headers(1).Range.Text = "Sravan Kumar Karpurapu"
